$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 33, shifting existing rows 33-117 down to 34-118
$ws.Rows.Item(33).Insert()

# Populate the new row 33 with the weekly Jengibre price record
$ws.Range("A33").Value = 8
$ws.Range("B33").Value = "Terminal La Palmera de La Serena"
$ws.Range("C33").Value = "Coquimbo"
$ws.Range("D33").Value = 45014
$ws.Range("E33").Value = 4
$ws.Range("F33").Value = 100114007
$ws.Range("G33").Value = "Jengibre"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 360
$ws.Range("K33").Value = 16000
$ws.Range("L33").Value = 17000
$ws.Range("M33").Value = 16500
$ws.Range("N33").Value = "$/caja 13 kilos"
$ws.Range("O33").Value = "Perú"
$ws.Range("P33").Value = 1269
$ws.Range("Q33").Value = 13
$ws.Range("R33").Value = "Hortaliza"
